$d = $word.ActiveDocument

# New (translated) date line that replaces every "V roku 2018 ..." sentence.
$newText = "V roku Cygnus: 10.-19. augusta, 9.-18. septembra, 8.-17. októbra"

# Walk every paragraph and replace any whose text starts with the old
# "V roku 2018 ..." sentence. Operating at the paragraph level (rather than
# just replacing the matched substring) ensures that when the sentence is
# split across multiple runs (constellation name run, trailing punctuation
# run, stray trailing-space run, etc.) the whole paragraph collapses down
# to a single, unformatted run - matching how Word behaves when you select
# the whole line and retype it.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "V roku 2018*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1   # exclude the paragraph mark
        $target = $d.Range($start, $end)
        $target.Text = ""
        $collapsed = $d.Range($start, $start)
        $collapsed.InsertAfter($newText)
    }
}
